$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.523.43'
$ws.Range("E2").Value = '  +0.21%  '
$ws.Range("D3").Value = '1.814.18'
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").Value = "'226.16"
$ws.Range("E5").Value = '  -0.88%  '
$ws.Range("E6").Value = '  +2.90%  '
$ws.Range("E7").Value = '  -0.21%  '
$ws.Range("D8").Value = "'38.09"
$ws.Range("E8").Value = '  +6.02%  '
$ws.Range("E9").Value = '  -3.96%  '
$ws.Range("E10").Value = '  -2.75%  '
$ws.Range("E11").Value = '  +1.12%  '
$ws.Range("D12").Value = '2.077.13'
$ws.Range("E12").Value = '  +0.24%  '
$ws.Range("D13").Value = "'11.26"
$ws.Range("E13").Value = '  -2.43%  '
$ws.Range("D14").Value = '1.813.23'
$ws.Range("E14").Value = '  -0.05%  '
$ws.Range("D15").Value = "'0.634"
$ws.Range("E15").Value = '  -2.05%  '
$ws.Range("D16").Value = '34.498.93'
$ws.Range("E16").Value = '  +0.12%  '
$ws.Range("E17").Value = '  -2.08%  '
$ws.Range("D18").Value = "'68.61"
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("D19").Value = "'243.87"
$ws.Range("E19").Value = '  -1.33%  '
$ws.Range("D20").Value = '0.0₃0776'
$ws.Range("D21").Value = "'11.33"
$ws.Range("E21").Value = '  -2.01%  '
$ws.Range("E22").Value = '  -0.23%  '
$ws.Range("D23").Value = "'4.14"
$ws.Range("E23").Value = '  -1.74%  '
$ws.Range("E24").Value = '  +3.61%  '
$ws.Range("D25").Value = "'170.38"
$ws.Range("E25").Value = '  -0.81%  '
$ws.Range("D26").Value = "'7.91"
$ws.Range("E26").Value = '  -1.03%  '
$ws.Range("E27").Value = '  +3.98%  '
$ws.Range("E28").Value = '  +0.71%  '
$ws.Range("E29").Value = '  -0.17%  '
$ws.Range("E30").Value = '  -0.89%  '
$ws.Range("E31").Value = '  -1.88%  '
$ws.Range("E32").Value = '  -2.54%  '
$ws.Range("E33").Value = '  -5.33%  '
$ws.Range("D34").Value = "'1.84"
$ws.Range("E34").Value = '  -0.26%  '
$ws.Range("D35").Value = '1.368.93'
$ws.Range("E35").Value = '  -2.34%  '
$ws.Range("D36").Value = "'0.649"
$ws.Range("E36").Value = '  -4.30%  '
$ws.Range("E37").Value = '  -0.84%  '
$ws.Range("D38").Value = "'2.35"
$ws.Range("E38").Value = '  -5.73%  '
$ws.Range("E39").Value = '  -1.98%  '
$ws.Range("D40").Value = "'1.22"
$ws.Range("E40").Value = '  -2.73%  '
$ws.Range("D41").Value = "'0.955"
$ws.Range("E41").Value = '  -1.35%  '
$ws.Range("E42").Value = '  +1.07%  '
$ws.Range("D43").Value = "'81.83"
$ws.Range("E43").Value = '  -1.56%  '
$ws.Range("D44").Value = "'2.79"
$ws.Range("E44").Value = '  -1.24%  '
$ws.Range("D45").Value = "'13.79"
$ws.Range("E45").Value = '  +2.79%  '
$ws.Range("E46").Value = '  +0.51%  '
$ws.Range("D47").Value = '1.977.74'
$ws.Range("E47").Value = '  +0.25%  '
$ws.Range("D48").Value = "'5.80"
$ws.Range("E48").Value = '  -4.23%  '
$ws.Range("E49").Value = '  -0.27%  '
$ws.Range("D50").Value = "'102.55"
$ws.Range("E50").Value = '  -3.04%  '
$ws.Range("E51").Value = '  -5.31%  '
